$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.055.40"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.893.65"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'307.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +2.01%  "
$ws.Range("D8").Value = "'0.3775"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").Value = "'0.07221"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'21.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "'0.8917"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "'0.07667"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "1.902.26"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "'94.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "'5.226"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008508"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'14.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "27.116.25"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'5.061"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "2.139.66"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").Value = "'10.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "'6.405"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'2.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.27%  "
$ws.Range("D26").Value = "'146.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").Value = "'1.734"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'114.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "'4.966"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.94%  "
$ws.Range("D31").Value = "'4.783"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "'0.09190"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("E34").Value = "  +6.17%  "
$ws.Range("D35").Value = "'0.7758"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("D36").Value = "'2.974"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'3.294"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("D38").Value = "'2.586"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "'0.5624"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "'0.01986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "'1.073"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'8.971"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").Value = "'6.624"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'118.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").Value = "'0.1516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").Value = "'0.4820"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "'1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "'10.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "'63.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
